# Generate Report for Handoff
# Adds a new handed-off file (b6e2c837-5f10-4add-93b4-b981bdb16f95.md) as a
# new row to the Overview / zh-cn / de-de worksheets, mirroring the existing
# row for 6c194b91-5ece-45e3-9832-a089f8429caf.md, and grows the three
# tables + autofilters + dimensions to include the new row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24b0b7a5163bb8cd02a23632c5b9f0f0c9ff10b6/e2e/b6e2c837-5f10-4add-93b4-b981bdb16f95.md"

# ---------------------------------------------------------------------
# Overview sheet -- new row 3
# ---------------------------------------------------------------------
$wsOverview.Range("A3").Value = "b6e2c837-5f10-4add-93b4-b981bdb16f95.md"

$wsOverview.Range("B3").Value = "e2e\b6e2c837-5f10-4add-93b4-b981bdb16f95.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseUrl, "", "", "e2e\b6e2c837-5f10-4add-93b4-b981bdb16f95.md") | Out-Null
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = 15570276

$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsOverview.Range("G3").Value = "2016-08-17 00:37:33"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# zh-cn sheet -- new row 3
# ---------------------------------------------------------------------
$wsZhCn.Range("A3").Value = "b6e2c837-5f10-4add-93b4-b981bdb16f95.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $baseUrl, "", "", "b6e2c837-5f10-4add-93b4-b981bdb16f95.md") | Out-Null
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = 15570276

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b6e2c837-5f10-4add-93b4-b981bdb16f95.9e5446f79b233590853bf287193639ff04ada25b.zh-cn.xlf"

$wsZhCn.Range("H3").Value = "2016-08-17 00:37:28"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# de-de sheet -- new row 3
# ---------------------------------------------------------------------
$wsDeDe.Range("A3").Value = "b6e2c837-5f10-4add-93b4-b981bdb16f95.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $baseUrl, "", "", "b6e2c837-5f10-4add-93b4-b981bdb16f95.md") | Out-Null
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = 15570276

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b6e2c837-5f10-4add-93b4-b981bdb16f95.9e5446f79b233590853bf287193639ff04ada25b.de-de.xlf"

$wsDeDe.Range("H3").Value = "2016-08-17 00:37:33"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
